# Text updates as supplied by PM&C.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsDesc = $wb.Worksheets.Item("Description")

# ---------------------------------------------------------------------------
# 1. Shared-string text corrections (Description sheet body copy)
# ---------------------------------------------------------------------------

# "...remained around 111 places..." -> "...remained around 108 places..."
$wsDesc.Range("B6").Value = "The growth in the rate of aged care places has stalled in recent years. After increasing in 2010 and 2011, rates of aged care places have remained around 108 places per 1000 people from 2012 to 2016."

# "...also influences aged care..." -> "...also influence aged care..."
$wsDesc.Range("B7").Value = "Australia" + [char]8217 + "s ageing population means that, even if service provision is maintained at similar levels, over time the number of aged care places per 1000 older people will decrease. Other factors, such as the availability of training for aged care professionals, also influence aged care availability. The Productivity Commission has also found that older Australians strongly prefer to age in place. Most people are happy staying in their family home, despite a common perception that such homes are too big for them (PC, 2015)."

# Notes text gains extra detail about Torres Strait Islander Australians / as-at date.
$wsDesc.Range("B8").Value = "Figures include operational number of aged care places per 1000 people aged 70 years or over and Aboriginal and Torres Strait Islander Australians aged 50-69 years as at 30 June."
$wsDesc.Range("B8").Font.Name = "Calibri"
$wsDesc.Range("B8").Font.Size = 11

# ---------------------------------------------------------------------------
# 2. Notes section restructure: old single "Sourced from" row becomes
#    "Source" + "References" rows with the new citation list.
# ---------------------------------------------------------------------------

$wsDesc.Range("A9").Value = "Source"
$wsDesc.Range("B9").Value = "Department of Health (unpublished)."

$wsDesc.Range("A10").Value = "References"
$wsDesc.Range("B10").Value = "Council of Australian Governments (COAG), 2011, National Healthcare Agreement."
$wsDesc.Range("B10").Font.Name = "Arial"
$wsDesc.Range("B10").Font.Size = 12
$wsDesc.Range("B10").WrapText = $true

$wsDesc.Range("B11").Value = "Department of Social Services (DSS), 2013, 2012-13 Report on the operation of the Aged Care Act 1997. Canberra: Department of Social Services."
$wsDesc.Range("B11").Font.Name = "Arial"
$wsDesc.Range("B11").Font.Size = 12
$wsDesc.Range("B11").WrapText = $true

$wsDesc.Range("B12").Value = "Productivity Commission (PC), 2015, Housing Decisions of Older Australians, Commission Research Paper, Canberra."
$wsDesc.Range("B12").Font.Name = "Arial"
$wsDesc.Range("B12").Font.Size = 12
$wsDesc.Range("B12").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Row heights (autofit-style adjustments that accompanied the rewrap)
# ---------------------------------------------------------------------------

$wsDesc.Rows.Item(6).RowHeight = 37.3
$wsDesc.Rows.Item(7).RowHeight = 73.1
$wsDesc.Rows.Item(8).RowHeight = 25.2
$wsDesc.Rows.Item(9).RowHeight = 13.8
$wsDesc.Rows.Item(10).RowHeight = 15
$wsDesc.Rows.Item(11).RowHeight = 27.25
$wsDesc.Rows.Item(12).RowHeight = 27.25

# ---------------------------------------------------------------------------
# 4. Selection / active-tab state: Description becomes the active sheet,
#    with B8:B12 (the edited notes block) selected there, while Data keeps
#    A1 as its (no-longer-active) selection.
# ---------------------------------------------------------------------------

$wsData.Range("A1").Select()
$wsDesc.Activate()
$wsDesc.Range("B8:B12").Select()
$wsDesc.Range("B12").Activate()
